# Update "想去人数" (want-to-go count) figures in the F column on the
# 展览 sheet and the matching rows on the 全部类型 sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets("展览")
$wsExhibit.Range("F2").Value = 62
$wsExhibit.Range("F3").Value = 11612
$wsExhibit.Range("F4").Value = 213
$wsExhibit.Range("F5").Value = 332
$wsExhibit.Range("F7").Value = 11583
$wsExhibit.Range("F10").Value = 86
$wsExhibit.Range("F11").Value = 1760
$wsExhibit.Range("F12").Value = 5737
$wsExhibit.Range("F13").Value = 115
$wsExhibit.Range("F14").Value = 3506
$wsExhibit.Range("F15").Value = 183
$wsExhibit.Range("F16").Value = 15

$wsAll = $wb.Worksheets("全部类型")
$wsAll.Range("F3").Value = 62
$wsAll.Range("F5").Value = 11612
$wsAll.Range("F6").Value = 213
$wsAll.Range("F7").Value = 332
$wsAll.Range("F9").Value = 11583
$wsAll.Range("F12").Value = 86
$wsAll.Range("F13").Value = 1760
$wsAll.Range("F15").Value = 5738
$wsAll.Range("F16").Value = 115
$wsAll.Range("F17").Value = 3506
$wsAll.Range("F18").Value = 183
$wsAll.Range("F19").Value = 15
